$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSE")

$data = @(
    @(5, "1.260567827915849E-4", "8.8609599328868442E-4"),
    @(6, "1.299822379837136E-4", "7.3875384176743753E-4"),
    @(7, "1.3078930621163439E-4", "7.0211228813735042E-4"),
    @(8, "1.2278006913106181E-4", "8.5502507398697804E-4"),
    @(9, "1.218396728429256E-4", "7.9082699613828536E-4"),
    @(10, "1.3126986203252259E-4", "8.0972213215268975E-4"),
    @(11, "1.288656729028824E-4", "5.93135595854183E-4"),
    @(12, "1.153423288888115E-4", "1.0602270190530651E-3"),
    @(13, "1.206247022007995E-4", "8.7573741621776134E-4"),
    @(14, "1.205017866461175E-4", "8.0682603068781456E-4"),
    @(15, "1.28494963953324E-4", "9.1650812249410016E-4"),
    @(16, "1.283765813976262E-4", "6.9479160147932354E-4"),
    @(17, "1.2234398612041269E-4", "7.2183592544594903E-4"),
    @(18, "1.211528045349201E-4", "1.1120615694756751E-3"),
    @(19, "1.2140786483355459E-4", "8.0205079880494472E-4"),
    @(20, "1.22644555519073E-4", "7.9212406907819646E-4"),
    @(21, "1.262607957238438E-4", "8.9358679777960348E-4"),
    @(22, "1.3167131078080751E-4", "6.997929800938936E-4"),
    @(23, "1.2734502024996179E-4", "9.2047436841062125E-4"),
    @(24, "1.264267009847704E-4", "9.4501680502124814E-4"),
    @(25, "1.208003758490001E-4", "1.031325876223435E-3"),
    @(26, "1.3090331775348111E-4", "6.9502658809301348E-4"),
    @(27, "1.19280870028474E-4", "1.008807098173565E-3"),
    @(28, "1.2649363776548889E-4", "8.0731188182762395E-4"),
    @(29, "1.2267562823929559E-4", "1.1892792231738541E-3"),
    @(30, "1.12157791070927E-4", "1.1369237280560169E-3"),
    @(31, "1.3044747241334231E-4", "6.6377997392327534E-4"),
    @(32, "1.1648007639381439E-4", "1.3587044815867361E-3"),
    @(33, "1.224841842411452E-4", "1.025608878342284E-3"),
    @(34, "1.272684550430182E-4", "8.3490364443467241E-4"),
    @(35, "1.2866750322555479E-4", "7.17364694464678E-4"),
    @(36, "1.3265605814514141E-4", "6.832482843200726E-4"),
    @(37, "1.2991186257052591E-4", "5.6686286723547073E-4"),
    @(38, "1.220250252783237E-4", "1.0504448754098451E-3"),
    @(39, "1.2209445548691709E-4", "9.3647869479658238E-4"),
    @(40, "1.2439359076342851E-4", "9.6545147770985571E-4"),
    @(41, "1.3096107737999469E-4", "7.6149643306337353E-4"),
    @(42, "1.2548107266803409E-4", "6.9647754537390427E-4"),
    @(43, "1.2740003126159101E-4", "9.3768594418841024E-4"),
    @(44, "1.330536623544695E-4", "8.8823340393249823E-4"),
    @(45, "1.2779039520215721E-4", "6.5136083260081341E-4"),
    @(46, "1.214142123615475E-4", "9.3091646188093481E-4"),
    @(47, "1.188251800856886E-4", "9.4420616755005787E-4"),
    @(48, "1.3009983434465161E-4", "9.093765616032235E-4"),
    @(49, "1.2297925447553749E-4", "1.116389392337704E-3"),
    @(50, "1.3396589627275889E-4", "8.4146355254407774E-4"),
    @(51, "1.2790272706201419E-4", "8.6374972545629402E-4"),
    @(52, "1.31336691587187E-4", "6.6900781955597766E-4"),
    @(53, "1.2736481373059661E-4", "8.5499840274208528E-4"),
    @(54, "1.208351987851741E-4", "9.7873719555829171E-4")
)

foreach ($item in $data) {
    $row = $item[0]
    $jval = [double]$item[1]
    $kval = [double]$item[2]
    $ws.Cells.Item($row, 10).Value = $jval
    $ws.Cells.Item($row, 11).Value = $kval
}

$ws.Select()
$excel.ActiveWindow.Zoom = 85
$ws.Range("M47").Select()

Write-Output "done"